# Add the new weekly-ranking worksheet for the 2025-11-19 snapshot,
# placed right after the most recent existing week's sheet.
$wb = $excel.ActiveWorkbook

$prevSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $prevSheet)
$newSheet.Name = "magapoke_2025-11-19"

# Copy the header cell formatting (bold font, thin border, centered
# alignment) from the previous week's sheet so the new sheet matches
# the existing style (same shared cellXf, no new styles introduced).
$prevSheet.Range("A1:B1").Copy()
$newSheet.Range("A1:B1").PasteSpecial(-4122)

$newSheet.Range("A1").Value = "rank"
$newSheet.Range("B1").Value = "title"

# Weekly ranking rows: rank (A) / title (B)
$newSheet.Cells.Item(2, 1).Value = 1
$newSheet.Cells.Item(2, 2).Value = '黒月のイェルクナハト'
$newSheet.Cells.Item(3, 1).Value = 2
$newSheet.Cells.Item(3, 2).Value = 'K-9~警視庁公安部公安第9課異能対策係~'
$newSheet.Cells.Item(4, 1).Value = 3
$newSheet.Cells.Item(4, 2).Value = 'アイドラトリィ'
$newSheet.Cells.Item(5, 1).Value = 4
$newSheet.Cells.Item(5, 2).Value = '黄昏町プリズナーズ'
$newSheet.Cells.Item(6, 1).Value = 5
$newSheet.Cells.Item(6, 2).Value = 'ゼロとヒャク'
$newSheet.Cells.Item(7, 1).Value = 6
$newSheet.Cells.Item(7, 2).Value = 'せいぶつ部の田辺くん'
$newSheet.Cells.Item(8, 1).Value = 7
$newSheet.Cells.Item(8, 2).Value = '篝家の８兄弟'
$newSheet.Cells.Item(9, 1).Value = 8
$newSheet.Cells.Item(9, 2).Value = 'ハードワーカー中田'
$newSheet.Cells.Item(10, 1).Value = 9
$newSheet.Cells.Item(10, 2).Value = 'MYS'
$newSheet.Cells.Item(11, 1).Value = 10
$newSheet.Cells.Item(11, 2).Value = '生きたがりの人狼'
$newSheet.Cells.Item(12, 1).Value = 11
$newSheet.Cells.Item(12, 2).Value = 'ともだちづくり'
$newSheet.Cells.Item(13, 1).Value = 12
$newSheet.Cells.Item(13, 2).Value = 'ナキナギ'
$newSheet.Cells.Item(14, 1).Value = 13
$newSheet.Cells.Item(14, 2).Value = 'ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜'
$newSheet.Cells.Item(15, 1).Value = 14
$newSheet.Cells.Item(15, 2).Value = '鉱石令嬢〜没落した悪役令嬢が炭鉱で一山当てるまでのお話〜'
$newSheet.Cells.Item(16, 1).Value = 15
$newSheet.Cells.Item(16, 2).Value = 'その青春'
$newSheet.Cells.Item(17, 1).Value = 16
$newSheet.Cells.Item(17, 2).Value = 'スルガメテオ'
$newSheet.Cells.Item(18, 1).Value = 17
$newSheet.Cells.Item(18, 2).Value = '夜鐘のキト'
$newSheet.Cells.Item(19, 1).Value = 18
$newSheet.Cells.Item(19, 2).Value = 'ドリーム☆ジャンボ☆ガール'
$newSheet.Cells.Item(20, 1).Value = 19
$newSheet.Cells.Item(20, 2).Value = 'お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！'
$newSheet.Cells.Item(21, 1).Value = 20
$newSheet.Cells.Item(21, 2).Value = '春くらり'
$newSheet.Cells.Item(22, 1).Value = 21
$newSheet.Cells.Item(22, 2).Value = '屋根の下のアルテミス'
$newSheet.Cells.Item(23, 1).Value = 22
$newSheet.Cells.Item(23, 2).Value = '皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～'
$newSheet.Cells.Item(24, 1).Value = 23
$newSheet.Cells.Item(24, 2).Value = 'ハナバス　苔石花江のバスケ論'
$newSheet.Cells.Item(25, 1).Value = 24
$newSheet.Cells.Item(25, 2).Value = 'それがメイドのカンナです'
$newSheet.Cells.Item(26, 1).Value = 25
$newSheet.Cells.Item(26, 2).Value = '卒業アルバムの彼女たち'
$newSheet.Cells.Item(27, 1).Value = 26
$newSheet.Cells.Item(27, 2).Value = '異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～'
$newSheet.Cells.Item(28, 1).Value = 27
$newSheet.Cells.Item(28, 2).Value = '白銀のキュイジーヌ～明治外交官の料理人～'
$newSheet.Cells.Item(29, 1).Value = 28
$newSheet.Cells.Item(29, 2).Value = '明智ナンバーワン'
$newSheet.Cells.Item(30, 1).Value = 29
$newSheet.Cells.Item(30, 2).Value = 'ナマイキ旭ちゃんをわからせたい'
$newSheet.Cells.Item(31, 1).Value = 30
$newSheet.Cells.Item(31, 2).Value = '追放されなかった男　～二度目の人生は土下座から始まりました～'
$newSheet.Cells.Item(32, 1).Value = 31
$newSheet.Cells.Item(32, 2).Value = '人生逆転ダンジョン'
$newSheet.Cells.Item(33, 1).Value = 32
$newSheet.Cells.Item(33, 2).Value = '平成転生'
$newSheet.Cells.Item(34, 1).Value = 33
$newSheet.Cells.Item(34, 2).Value = '限界集落を脱村した錬金術士、都会で"最強"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～'
$newSheet.Cells.Item(35, 1).Value = 34
$newSheet.Cells.Item(35, 2).Value = '鳴るさんだぁ'
$newSheet.Cells.Item(36, 1).Value = 35
$newSheet.Cells.Item(36, 2).Value = '永久のユウグレ'
$newSheet.Cells.Item(37, 1).Value = 36
$newSheet.Cells.Item(37, 2).Value = '花子狩り'
$newSheet.Cells.Item(38, 1).Value = 37
$newSheet.Cells.Item(38, 2).Value = 'JK Biker'
$newSheet.Cells.Item(39, 1).Value = 38
$newSheet.Cells.Item(39, 2).Value = 'ハプスブルク家の華麗なる受難'
$newSheet.Cells.Item(40, 1).Value = 39
$newSheet.Cells.Item(40, 2).Value = '〈小市民〉 春期限定いちごタルト事件'
$newSheet.Cells.Item(41, 1).Value = 40
$newSheet.Cells.Item(41, 2).Value = 'じゅーくぼっくす'
$newSheet.Cells.Item(42, 1).Value = 41
$newSheet.Cells.Item(42, 2).Value = '眠れる森のレガ'
$newSheet.Cells.Item(43, 1).Value = 42
$newSheet.Cells.Item(43, 2).Value = '東京デスレース'
$newSheet.Cells.Item(44, 1).Value = 43
$newSheet.Cells.Item(44, 2).Value = 'イエティ、とある日々'

# Restore the originally active sheet/tab
$wb.Worksheets.Item(1).Activate()
